$d = $word.ActiveDocument

# Helper: split the run containing character position $pos into two runs at
# that exact position, without leaving any formatting residue. Adding and
# then immediately deleting a bookmark at a collapsed point forces the host
# to re-materialize the run boundary there cleanly (no stray <w:rPr/>).
function SplitAt($pos) {
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add("tmp_split_marker", $r) | Out-Null
    $d.Bookmarks("tmp_split_marker").Delete()
}

# ---------------------------------------------------------------------------
# Fix 1 (Step 15 paragraph): remove a duplicated left smart-quote typo:
#   under ““Plant,   ->   under “Plant,
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "under " + [char]0x201C + [char]0x201C + "Plant, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "under " + [char]0x201C + "Plant, ", 2
) | Out-Null

# ---------------------------------------------------------------------------
# Fix 2 (Step 16 paragraph): add the missing file-name reference.
#   : Locate the Weightings H:A (animal Husbandry: Agriculture) table ...
#   ->
#   : Locate the "relative_weightings_H-A.csv" (animal Husbandry: Agriculture) table ...
# and split the edited run back up into the 5 runs shown by the diff, while
# restoring every other pre-existing run boundary in that paragraph (a text
# length change elsewhere in the same paragraph can otherwise cause the
# engine to coalesce neighbouring same-formatted runs).
# ---------------------------------------------------------------------------

$openQuote  = [string][char]0x201C
$closeQuote = [string][char]0x201D
$fileRef    = $openQuote + "relative_weightings_H-A.csv" + $closeQuote

$find1 = $d.Content
$find1.Find.Execute("Weightings H:A") | Out-Null
$wStart = $find1.Start
$wEnd   = $find1.End

$target = $d.Range($wStart, $wEnd)
$target.Text = $fileRef

# Re-locate (fresh, post-edit) each subsequent original run's starting text
# so the paragraph's pre-existing run boundaries (which may have been
# coalesced by the edit above) get restored exactly as they were.
$anchors = @(
    " (animal Husbandry)",
    " ",
    "weighting ",
    "value for each society and every value related to the plant environmental data, by the " + $openQuote + "A" + $closeQuote,
    " (Agriculture)",
    " ",
    "weighting ",
    "value for each society."
)

$searchFrom = $wStart
foreach ($anchor in $anchors) {
    $r = $d.Range($searchFrom, $d.Content.End)
    $r.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    SplitAt $r.Start
    $searchFrom = $r.Start + $anchor.Length
}

# Now split the edited run itself into the 5 pieces shown by the diff:
#   ": Locate the "  |  """  |  "relative_weightings_H-A"  |  ".csv"""  |  " (animal Husbandry: ..."
$nameStart = $wStart + $openQuote.Length
$nameEnd   = $nameStart + "relative_weightings_H-A".Length
$refEnd    = $nameEnd + (".csv" + $closeQuote).Length

SplitAt $wStart
SplitAt $nameStart
SplitAt $nameEnd
SplitAt $refEnd
